# Update "models performance metrics.xlsx"
# Adds two new summary tables to the "Testing Metrics" sheet:
#   - "Testing trained models on our dataset"  (rows 25-30)
#   - "Training the models on our dataset"     (rows 34-39)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testing Metrics")

# ---------------------------------------------------------------------
# Section 1: "Testing trained models on our dataset" (rows 25-30)
# ---------------------------------------------------------------------
$ws.Range("B25").Value = "Testing trained models on our dataset"

$ws.Range("F26").Value = "Accuracy"
$ws.Range("G26").Value = "Loss"

$ws.Range("B27").Value = "ANN"
$ws.Range("B28").Value = "CNN Spectrogram"
$ws.Range("B29").Value = "1-D CNN with Flattened Spectrogram"
$ws.Range("B30").Value = "1-D CNN with MFCC"

# Carry over the number formatting/font used by the other F/G data columns
$ws.Range("F11:G11").Copy()
$ws.Range("F27:G30").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F27").Value = 0.66279071569442705
$ws.Range("G27").Value = 9.7774999999999999

$ws.Range("F28").Value = 0.41860464215278598
$ws.Range("G28").Value = 4.9439000000000002

$ws.Range("F29").Value = 0.11627907305955799
$ws.Range("G29").Value = 2.0794999999999999

$ws.Range("F30").Value = 0.62790697813034002
$ws.Range("G30").Value = 5.1040999999999999

# ---------------------------------------------------------------------
# Section 2: "Training the models on our dataset" (rows 34-39)
# ---------------------------------------------------------------------
$ws.Range("B34").Value = "Training the models on our dataset"

$ws.Range("F35").Value = "Accuracy"
$ws.Range("G35").Value = "Loss"
$ws.Range("H35").Value = "epochs"

$ws.Range("B36").Value = "ANN"
$ws.Range("B37").Value = "CNN Spectrogram"
$ws.Range("B38").Value = "1-D CNN with Flattened Spectrogram"
$ws.Range("B39").Value = "1-D CNN with MFCC"

$ws.Range("F11:G11").Copy()
$ws.Range("F36:G39").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("H8").Copy()
$ws.Range("H36:H39").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F36").Value = 0.1201
$ws.Range("G36").Value = 2.0933000000000002
$ws.Range("H36").Value = 200

$ws.Range("F37").Value = 0.98970000000000002
$ws.Range("G37").Value = 0.0089
$ws.Range("H37").Value = 10

$ws.Range("F38").Value = 0.1298
$ws.Range("G38").Value = 2.0783
$ws.Range("H38").Value = 20

$ws.Range("F39").Value = 0.89534884691238403
$ws.Range("G39").Value = 0.16669999999999999
$ws.Range("H39").Value = 30

# ---------------------------------------------------------------------
# View state: scroll / zoom / selection on the "Testing Metrics" sheet
# ---------------------------------------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$win.Zoom = 115
$ws.Range("G39").Select()
